# Re-sort the "Periodo Mora" column (ascending instead of descending) for both
# employees' detail rows, and update the associated Valor Mora / Salario Basico
# figures for the second employee's block, matching the refreshed account
# statement data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First employee block: rows 16-33 (Periodo Mora 1609..1802, ascending) ---
$periodos1 = @("1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802")

for ($i = 0; $i -lt $periodos1.Count; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos1[$i]
}

# --- Second employee block: rows 34-43 (Periodo Mora 2005..2102, ascending) ---
$periodos2 = @("2005","2006","2007","2008","2009","2010","2011","2012","2101","2102")

for ($i = 0; $i -lt $periodos2.Count; $i++) {
    $row = 34 + $i
    $ws.Range("E$row").Value = $periodos2[$i]
}

# Updated "Salario Basico" (col F) and "Valor Mora" (col G) for the second
# employee's block: all periods now carry 48000/1200000, except the last
# period (2102) which keeps the odd Salario Basico value of 35200.
for ($row = 34; $row -le 42; $row++) {
    $ws.Range("F$row").Value = 48000
    $ws.Range("G$row").Value = 1200000
}
$ws.Range("F43").Value = 35200
$ws.Range("G43").Value = 1200000
